$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Original layout (row 1 = header "Conta/Nome/Saldo"):
#   row 3: 004204344 / CLINEO    / 49456.58
#   row 4: 004211807 / EDINARDO  / 22089.92
#   row 5: 005690206 / KAUANNE   / 20020
#   row 6: 005000460 / MARIANA   / 4500
#
# Target layout:
#   row 3: 005206566 / LEVI      / 40000      (replaces the CLINEO+EDINARDO pair)
#   row 4: 005690206 / KAUANNE   / 20020      (unchanged, shifts up)
#   (MARIANA row removed entirely)

# Remove the EDINARDO row (row 4) so CLINEO's row (now row 3) can be
# overwritten with LEVI's data, collapsing the two rows into one.
$ws.Rows.Item(4).Delete()

# Overwrite the remaining row (originally CLINEO) with the new record.
# The account number has a leading zero, so force text formatting while
# assigning it (otherwise it is auto-coerced to a number and the zero is
# lost), then clear the temporary format so the cell's style matches the
# plain (unstyled) data cells around it.
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "005206566"
$ws.Range("A3").ClearFormats()

$ws.Range("B3").Value = "LEVI"
$ws.Range("C3").Value = 40000

# Remove the MARIANA row, which is now row 5 after the earlier deletion.
$ws.Rows.Item(5).Delete()
